# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook stores this date as an Excel serial number; the automatic
# update increments it from 46060 to 46061 (one day later) for every
# data row (rows 2 through 85).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row on the sheet (data rows start at row 2,
# row 1 holds the headers). Fall back to the known last row (85) if the
# used range can't be determined for some reason.
$lastRow = 85
$used = $ws.UsedRange
if ($used -ne $null) {
    $computedLast = $used.Row + $used.Rows.Count - 1
    if ($computedLast -ge 2) {
        $lastRow = $computedLast
    }
}

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = 46061
}
